$p = $ppt.ActivePresentation

# --- Slide 1: title slide subtitle "Haiyue Wang" -> blank paragraph + "Haiyue Wang" split into runs ---
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$subTr = $subtitle.TextFrame.TextRange
[void]$subTr.InsertBefore("`r")
$subTr2 = $subtitle.TextFrame.TextRange
$namePara = $subTr2.Paragraphs(2, 1)
$nameRun1 = $namePara.Characters(1, 6)
$nameRun1.Text = "Haiyue"
$nameRun2 = $namePara.Characters(7, 1)
$nameRun2.Text = " "
$nameRun3 = $namePara.Characters(8, 4)
$nameRun3.Text = "Wang"

# --- Slide 2: Outline bullets 2 & 3 reworded ---
$s2 = $p.Slides.Item(2)
$outline = $s2.Shapes.Item(2)
$outlineTr = $outline.TextFrame.TextRange

$bullet2 = $outlineTr.Paragraphs(2, 1)
$bullet2.Text = "2. The Domain using Kalma Filter"
$b2r1 = $bullet2.Characters(1, 7)
$b2r1.Text = "2. The "
$b2r2 = $bullet2.Characters(8, 13)
$b2r2.Text = "Domain using "
$b2r3 = $bullet2.Characters(21, 5)
$b2r3.Text = "Kalma"
$b2r4 = $bullet2.Characters(26, 7)
$b2r4.Text = " Filter"

$bullet3 = $outlineTr.Paragraphs(3, 1)
$bullet3.Text = "3. Three examples"
$b3r1 = $bullet3.Characters(1, 3)
$b3r1.Text = "3. "
$b3r2 = $bullet3.Characters(4, 14)
$b3r2.Text = "Three examples"

# --- Slide 3: highlight key words in red within the definition paragraph ---
$s3 = $p.Slides.Item(3)
$body = $s3.Shapes.Item(2)
$bodyTr = $body.TextFrame.TextRange
$defPara = $bodyTr.Paragraphs(1, 1)

$d1 = $defPara.Characters(1, 55)
$d1.Text = "The Kalman Filter is a mathematical algorithm used for "

$d2 = $defPara.Characters(56, 10)
$d2.Text = "estimating"
$d2.Font.Color.RGB = 255

$d3 = $defPara.Characters(66, 5)
$d3.Text = " and "

$d4 = $defPara.Characters(71, 10)
$d4.Text = "predicting"
$d4.Font.Color.RGB = 255

$d5 = $defPara.Characters(81, 16)
$d5.Text = " the state of a "

$d6 = $defPara.Characters(97, 7)
$d6.Text = "dynamic"
$d6.Font.Color.RGB = 255

$d7 = $defPara.Characters(104, 41)
$d7.Text = " system, particularly in the presence of "

$d8 = $defPara.Characters(145, 23)
$d8.Text = "noisy or uncertain data"
$d8.Font.Color.RGB = 255

$d9 = $defPara.Characters(168, 1)
$d9.Text = "."
